$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Range("F2").Value = 27.16937875747681
$ws.Range("F3").Value = 26.72026777267456
$ws.Range("F4").Value = 26.70567917823792
$ws.Range("F5").Value = 26.76639604568481
$ws.Range("F6").Value = 26.73976349830627
$ws.Range("F7").Value = 26.7776243686676
$ws.Range("F8").Value = 26.66259741783142
$ws.Range("F9").Value = 26.97221636772156
$ws.Range("F10").Value = 26.62291574478149
$ws.Range("F11").Value = 26.94012475013733
$ws.Range("F12").Value = 26.57977604866028
$ws.Range("F13").Value = 26.67786335945129
$ws.Range("F14").Value = 26.7407374382019
$ws.Range("F15").Value = 26.71946358680725
$ws.Range("F16").Value = 26.66970825195312
$ws.Range("F17").Value = 26.83515858650208
$ws.Range("F18").Value = 26.61966824531555
$ws.Range("F19").Value = 26.75222730636597
$ws.Range("F20").Value = 26.7217960357666
$ws.Range("F21").Value = 26.78909873962402

$ws = $wb.Worksheets.Item("run_2")
$ws.Range("F2").Value = 26.8798987865448
$ws.Range("F3").Value = 26.75906491279602
$ws.Range("F4").Value = 26.79637861251831
$ws.Range("F5").Value = 26.56688904762268
$ws.Range("F6").Value = 26.82339334487915
$ws.Range("F7").Value = 26.70717406272888
$ws.Range("F8").Value = 26.60329842567444
$ws.Range("F9").Value = 26.64246559143066
$ws.Range("F10").Value = 26.66719722747803
$ws.Range("F11").Value = 26.85017108917236
$ws.Range("F12").Value = 26.79422760009766
$ws.Range("F13").Value = 26.64251446723938
$ws.Range("F14").Value = 26.67646050453186
$ws.Range("F15").Value = 26.55833554267884
$ws.Range("F16").Value = 26.65101337432861
$ws.Range("F17").Value = 26.60952568054199
$ws.Range("F18").Value = 26.56052279472351
$ws.Range("F19").Value = 26.68291783332825
$ws.Range("F20").Value = 26.63772368431092
$ws.Range("F21").Value = 26.87133479118347

$ws = $wb.Worksheets.Item("run_3")
$ws.Range("F2").Value = 26.89574003219604
$ws.Range("F3").Value = 26.67002391815185
$ws.Range("F4").Value = 26.637540102005
$ws.Range("F5").Value = 26.66706418991089
$ws.Range("F6").Value = 26.62338161468506
$ws.Range("F7").Value = 26.80319428443909
$ws.Range("F8").Value = 26.69219708442688
$ws.Range("F9").Value = 26.71279788017273
$ws.Range("F10").Value = 26.79429697990417
$ws.Range("F11").Value = 26.91580033302307
$ws.Range("F12").Value = 26.68517470359802
$ws.Range("F13").Value = 26.71406579017639
$ws.Range("F14").Value = 26.75562930107117
$ws.Range("F15").Value = 26.60142612457276
$ws.Range("F16").Value = 26.59487366676331
$ws.Range("F17").Value = 26.65093779563904
$ws.Range("F18").Value = 26.64496350288391
$ws.Range("F19").Value = 26.65395545959473
$ws.Range("F20").Value = 26.58806443214417
$ws.Range("F21").Value = 26.96683692932129

$ws = $wb.Worksheets.Item("run_4")
$ws.Range("F2").Value = 26.977623462677
$ws.Range("F3").Value = 26.73638725280762
$ws.Range("F4").Value = 26.67285299301147
$ws.Range("F5").Value = 26.58229804039001
$ws.Range("F6").Value = 26.65340495109558
$ws.Range("F7").Value = 26.68390989303589
$ws.Range("F8").Value = 26.6735634803772
$ws.Range("F9").Value = 26.70522975921631
$ws.Range("F10").Value = 26.66165399551392
$ws.Range("F11").Value = 26.74085474014282
$ws.Range("F12").Value = 26.54023122787476
$ws.Range("F13").Value = 26.571861743927
$ws.Range("F14").Value = 26.45703983306885
$ws.Range("F15").Value = 26.59520936012268
$ws.Range("F16").Value = 26.47363662719727
$ws.Range("F17").Value = 26.7348153591156
$ws.Range("F18").Value = 26.46002411842347
$ws.Range("F19").Value = 26.61373233795166
$ws.Range("F20").Value = 26.501291513443
$ws.Range("F21").Value = 26.89356255531311

$ws = $wb.Worksheets.Item("run_5")
$ws.Range("F2").Value = 26.80537438392639
$ws.Range("F3").Value = 26.61231660842896
$ws.Range("F4").Value = 26.70644783973694
$ws.Range("F5").Value = 27.66968560218811
$ws.Range("F6").Value = 26.59879875183105
$ws.Range("F7").Value = 26.51203632354736
$ws.Range("F8").Value = 26.456209897995
$ws.Range("F9").Value = 26.47403597831726
$ws.Range("F10").Value = 26.49908638000488
$ws.Range("F11").Value = 26.67487335205078
$ws.Range("F12").Value = 26.57961916923523
$ws.Range("F13").Value = 26.71573758125305
$ws.Range("F14").Value = 26.68586444854736
$ws.Range("F15").Value = 26.79987096786499
$ws.Range("F16").Value = 27.34831547737122
$ws.Range("F17").Value = 26.91782093048096
$ws.Range("F18").Value = 26.97077989578247
$ws.Range("F19").Value = 26.83722329139709
$ws.Range("F20").Value = 26.89557528495789
$ws.Range("F21").Value = 27.19285130500793
